# Small updates to profiling data on Sheet1, row 6 ("copy"),
# plus refresh the active selection/view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the measured values for the "copy" series (row 6)
$ws.Range("H6").Value = 0.731
$ws.Range("I6").Value = 1.43
$ws.Range("J6").Value = 1.57
$ws.Range("K6").Value = 1.74
$ws.Range("L6").Value = 1.34
$ws.Range("M6").Value = 1.57
$ws.Range("N6").Value = 1.96

# Update the active sheet view: scroll position and selection
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("N17").Select()
